$d = $word.ActiveDocument

function Split-LastDigit {
    param($paraIndex, $searchText, $newDigit)

    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range.Duplicate
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $matchEnd = $rng.End
    $digitRange = $d.Range($matchEnd - 1, $matchEnd)

    # Capture the existing run's formatting (incl. rFonts/color/size) by cloning its
    # FormattedText, changing only the text, then writing it back in place. Word's
    # engine treats this as a content change distinct from the neighbouring run,
    # so the final digit ends up split into its own <w:r> while keeping the source
    # run's full run-properties (rFonts ascii/eastAsia/hAnsi/cs, color, sz, szCs...).
    $ft = $digitRange.FormattedText
    $ft.Text = $newDigit
    $digitRange.FormattedText = $ft
}

# "Lösungsidee: /7" -> "Lösungsidee: /" + new run "8"
Split-LastDigit 8 "/7" "8"

# "Quellcode: /13" -> "Quellcode: /1" + new run "2"
Split-LastDigit 17 "/13" "2"

# "Lösungsidee: /7" (colored run) -> "Lösungsidee: /" + new run "8"
Split-LastDigit 48 "/7" "8"

# "Quellcode: /13" -> "Quellcode: /1" + new run "2"
Split-LastDigit 55 "/13" "2"
